$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.438353082513629
$ws.Range("C2").Value = 0.02551663813159166
$ws.Range("D2").Value = 0.6301285277367725
$ws.Range("E2").Value = 0.2434783830575107
$ws.Range("G2").Value = 0.002533955021152286
$ws.Range("J2").Value = 0.1165594566598962
$ws.Range("K2").Value = 0.3861236454256414
$ws.Range("M2").Value = 0.3240210665798813
$ws.Range("O2").Value = 6.493975059679684

$ws.Range("B3").Value = 0.4099427871376804
$ws.Range("C3").Value = 0.02238307421842478
$ws.Range("D3").Value = 0.6248026463177325
$ws.Range("E3").Value = 0.2420530726684511
$ws.Range("G3").Value = 0.002537101871097631
$ws.Range("J3").Value = 0.1163041069418149
$ws.Range("K3").Value = 0.3571244798797863
$ws.Range("M3").Value = 0.3143606908792478
$ws.Range("O3").Value = 6.499767941423812

$ws.Range("B4").Value = 0.3927038230266646
$ws.Range("C4").Value = 0.02045940630277698
$ws.Range("D4").Value = 0.6218220577688243
$ws.Range("E4").Value = 0.2412857709669431
$ws.Range("G4").Value = 0.002539137413494199
$ws.Range("J4").Value = 0.1161953437903165
$ws.Range("K4").Value = 0.3394821284449421
$ws.Range("M4").Value = 0.3085963874768112
$ws.Range("O4").Value = 6.506022152374243

$ws.Range("B5").Value = 0.3857306885595051
$ws.Range("C5").Value = 0.0196756155546467
$ws.Range("D5").Value = 0.6206803455029473
$ws.Range("E5").Value = 0.2410002372427833
$ws.Range("G5").Value = 0.00253999298485402
$ws.Range("J5").Value = 0.1161631061015846
$ws.Range("K5").Value = 0.3323340913769073
$ws.Range("M5").Value = 0.3062895471690439
$ws.Range("O5").Value = 6.509249082385139

$ws.Range("B6").Value = 0.3845759479010553
$ws.Range("C6").Value = 0.01954547583486743
$ws.Range("D6").Value = 0.6204951710610374
$ws.Range("E6").Value = 0.2409544651316722
$ws.Range("G6").Value = 0.002540136628855498
$ws.Range("J6").Value = 0.1161584831660036
$ws.Range("K6").Value = 0.331149672861855
$ws.Range("M6").Value = 0.3059090478663791
$ws.Range("O6").Value = 6.509825879186963

$ws.Range("B7").Value = 0.3926095704012482
$ws.Range("C7").Value = 0.02044883530173536
$ws.Range("D7").Value = 0.6218063649296539
$ws.Range("E7").Value = 0.2412818101987213
$ws.Range("G7").Value = 0.002539148846319973
$ws.Range("J7").Value = 0.116194860081599
$ws.Range("K7").Value = 0.3393855596172131
$ws.Range("M7").Value = 0.3085651057636554
$ws.Range("O7").Value = 6.506062925515465

$ws.Range("B8").Value = 0.4285148292560166
$ws.Range("C8").Value = 0.02443612940361106
$ws.Range("D8").Value = 0.6282321288906729
$ws.Range("E8").Value = 0.242964569537893
$ws.Range("G8").Value = 0.002535018649774666
$ws.Range("J8").Value = 0.1164614497278258
$ws.Range("K8").Value = 0.3760910271105899
$ws.Range("M8").Value = 0.3206555301103009
$ws.Range("O8").Value = 6.495412604876066

$ws.Range("B9").Value = 0.5005422801220618
$ws.Range("C9").Value = 0.03225705669143508
$ws.Range("D9").Value = 0.6431273292498645
$ws.Range("E9").Value = 0.2471193035042276
$ws.Range("G9").Value = 0.002527735853214345
$ws.Range("J9").Value = 0.117365048628983
$ws.Range("K9").Value = 0.4493560112061914
$ws.Range("M9").Value = 0.3456881148516473
$ws.Range("O9").Value = 6.49593508864865

$ws.Range("B10").Value = 0.5544396856747653
$ws.Range("C10").Value = 0.03800355429557101
$ws.Range("D10").Value = 0.6554675828647589
$ws.Range("E10").Value = 0.2506924433915003
$ws.Range("G10").Value = 0.002522877845269993
$ws.Range("J10").Value = 0.1182610102285437
$ws.Range("K10").Value = 0.5039606093600071
$ws.Range("M10").Value = 0.3648841080317737
$ws.Range("O10").Value = 6.509386009150148

$ws.Range("B11").Value = 0.5791704531911535
$ws.Range("C11").Value = 0.04061778737623456
$ws.Range("D11").Value = 0.661384513115479
$ws.Range("E11").Value = 0.2524309585931448
$ws.Range("G11").Value = 0.002520773703027862
$ws.Range("J11").Value = 0.1187190046468842
$ws.Range("K11").Value = 0.5289692711807845
$ws.Range("M11").Value = 0.3737912711999414
$ws.Range("O11").Value = 6.518346742162805

$ws.Range("B12").Value = 0.5885656950268583
$ws.Range("C12").Value = 0.04160772614774544
$ws.Range("D12").Value = 0.6636686499356586
$ws.Range("E12").Value = 0.2531055317195694
$ws.Range("G12").Value = 0.002519992049385
$ws.Range("J12").Value = 0.1188996811642511
$ws.Range("K12").Value = 0.5384634486333368
$ws.Range("M12").Value = 0.3771892440147937
$ws.Range("O12").Value = 6.522148783708417

$ws.Range("B13").Value = 0.5865409220318156
$ws.Range("C13").Value = 0.04139452627701701
$ws.Range("D13").Value = 0.6631747858773451
$ws.Range("E13").Value = 0.2529595286608952
$ws.Range("G13").Value = 0.002520159720484915
$ws.Range("J13").Value = 0.1188604471593635
$ws.Range("K13").Value = 0.5364176480870242
$ws.Range("M13").Value = 0.3764563186534815
$ws.Range("O13").Value = 6.52131176179131

$ws.Range("B14").Value = 0.5799428012809642
$ws.Range("C14").Value = 0.04069923070910875
$ws.Range("D14").Value = 0.6615715586516728
$ws.Range("E14").Value = 0.2524861308961164
$ws.Range("G14").Value = 0.002520709092877057
$ws.Range("J14").Value = 0.1187337238659367
$ws.Range("K14").Value = 0.5297498849443798
$ws.Range("M14").Value = 0.3740703236113347
$ws.Range("O14").Value = 6.518651344435511

$ws.Range("B15").Value = 0.5759051897863117
$ws.Range("C15").Value = 0.04027333925544951
$ws.Range("D15").Value = 0.6605952011408078
$ws.Range("E15").Value = 0.25219827465272
$ws.Range("G15").Value = 0.002521047568963842
$ws.Range("J15").Value = 0.1186570453943361
$ws.Range("K15").Value = 0.5256687955259167
$ws.Range("M15").Value = 0.3726120881756145
$ws.Range("O15").Value = 6.517075005231732

$ws.Range("B16").Value = 0.5528277209913028
$ws.Range("C16").Value = 0.03783270771924663
$ws.Range("D16").Value = 0.6550869928848329
$ws.Range("E16").Value = 0.2505811005669827
$ws.Range("G16").Value = 0.002523017478478849
$ws.Range("J16").Value = 0.1182320932257355
$ws.Range("K16").Value = 0.5023296022907004
$ws.Range("M16").Value = 0.3643055118247958
$ws.Range("O16").Value = 6.508857602322507

$ws.Range("B17").Value = 0.5387246494320834
$ws.Range("C17").Value = 0.03633546766087647
$ws.Range("D17").Value = 0.6517855014190843
$ws.Range("E17").Value = 0.2496179602308928
$ws.Range("G17").Value = 0.002524252997209306
$ws.Range("J17").Value = 0.1179843075845142
$ws.Range("K17").Value = 0.4880547522235759
$ws.Range("M17").Value = 0.3592543878504912
$ws.Range("O17").Value = 6.504544463828495

$ws.Range("B18").Value = 0.5306329631185065
$ws.Range("C18").Value = 0.03547430820380271
$ws.Range("D18").Value = 0.6499151229186282
$ws.Range("E18").Value = 0.2490746321117427
$ws.Range("G18").Value = 0.002524973596768554
$ws.Range("J18").Value = 0.1178465336723065
$ws.Range("K18").Value = 0.479860139532974
$ws.Range("M18").Value = 0.3563655756052171
$ws.Range("O18").Value = 6.502331138881658

$ws.Range("B19").Value = 0.5278967069922942
$ws.Range("C19").Value = 0.03518273770181679
$ws.Range("D19").Value = 0.649286751447022
$ws.Range("E19").Value = 0.2488924995544863
$ws.Range("G19").Value = 0.002525219292517522
$ws.Range("J19").Value = 0.1178007010769733
$ws.Range("K19").Value = 0.4770883259289178
$ws.Range("M19").Value = 0.3553903040373783
$ws.Range("O19").Value = 6.501627680132344

$ws.Range("B20").Value = 0.5402238751096036
$ws.Range("C20").Value = 0.03649485026615196
$ws.Range("D20").Value = 0.652133996433804
$ws.Range("E20").Value = 0.2497193866568068
$ws.Range("G20").Value = 0.002524120443624343
$ws.Range("J20").Value = 0.1180101936561755
$ws.Range("K20").Value = 0.4895726908065114
$ws.Range("M20").Value = 0.3597903858321203
$ws.Range("O20").Value = 6.504975920669892

$ws.Range("B21").Value = 0.5818800112846816
$ws.Range("C21").Value = 0.04090345652069516
$ws.Range("D21").Value = 0.6620412848307637
$ws.Range("E21").Value = 0.2526247388762641
$ws.Range("G21").Value = 0.002520547318612603
$ws.Range("J21").Value = 0.1187707490165835
$ws.Range("K21").Value = 0.5317077214764367
$ws.Range("M21").Value = 0.3747704699679844
$ws.Range("O21").Value = 6.519421677380706

$ws.Range("B22").Value = 0.6092807966137457
$ws.Range("C22").Value = 0.04378463609161543
$ws.Range("D22").Value = 0.6687699195455821
$ws.Range("E22").Value = 0.2546181739167324
$ws.Range("G22").Value = 0.002518300285260044
$ws.Range("J22").Value = 0.1193100366401367
$ws.Range("K22").Value = 0.5593847904614506
$ws.Range("M22").Value = 0.3847066243547772
$ws.Range("O22").Value = 6.531245723973086

$ws.Range("B23").Value = 0.5946404777999419
$ws.Range("C23").Value = 0.04224691594363605
$ws.Range("D23").Value = 0.6651555407711669
$ws.Range("E23").Value = 0.2535455906528696
$ws.Range("G23").Value = 0.002519491522034461
$ws.Range("J23").Value = 0.1190183474193631
$ws.Range("K23").Value = 0.5446003709361946
$ws.Range("M23").Value = 0.3793902081721541
$ws.Range("O23").Value = 6.524716920227775

$ws.Range("B24").Value = 0.539546024206544
$ws.Range("C24").Value = 0.036422794560238
$ws.Range("D24").Value = 0.6519763555728844
$ws.Range("E24").Value = 0.2496734993917968
$ws.Range("G24").Value = 0.00252418033914326
$ws.Range("J24").Value = 0.1179984759811035
$ws.Range("K24").Value = 0.4888863927836269
$ws.Range("M24").Value = 0.3595480139582321
$ws.Range("O24").Value = 6.504780029329623

$ws.Range("B25").Value = 0.4808845454911364
$ws.Range("C25").Value = 0.03014116089734387
$ws.Range("D25").Value = 0.6388523629105691
$ws.Range("E25").Value = 0.2459038711254884
$ws.Range("G25").Value = 0.002529619159613223
$ws.Range("J25").Value = 0.1170798395845907
$ws.Range("K25").Value = 0.4293989825189897
$ws.Range("M25").Value = 0.3387746964112637
$ws.Range("O25").Value = 6.49350022113282
